# Incorporate new urban net returns data: add a "fipsstr" variable row into
# the pointpanel_variables data dictionary (inserted right after "fips"),
# which pushes every subsequent row down by one. Also refresh the urban_nr
# variable label text, and re-sequence the "position" index column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 8 (shifts old rows 8..58 down to 9..59,
#    which also naturally reproduces the duplicated last "lccL78_acresk"
#    row at the new row 59).
$ws.Rows.Item(8).Insert()

# 2) Populate the newly inserted row 8 with the new "fipsstr" variable.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "fipsstr"
$ws.Cells.Item(8, 3).Value = "str5"
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = ""

# 3) Re-sequence the "position" column (A) for every shifted row so it
#    again reads 0..58 top to bottom (A holds a literal number, not a
#    formula, so the insert doesn't renumber it automatically).
for ($r = 9; $r -le 59; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# 4) Update the urban_nr variable label (row 33 after the shift) to the
#    new RFF-sourced description.
$ws.Cells.Item(33, 5).Value = "2010USD annualized net return/acre [RFF]"
